$d = $word.ActiveDocument

# --- Locate the paragraph that contains the typo ("ehat" -> "chat") ---
# Find the start of the paragraph text.
$startRange = $d.Content
$startRange.Find.Execute("Divided up the main tasks", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rStart = $startRange.Start

# Find the end of the paragraph text (the final sentence).
$endRange = $d.Content
$endRange.Find.Execute("they are completed.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rEnd = $endRange.End

# The three pieces of corrected text (the middle run ends right after the
# fixed character "c", matching where the _GoBack bookmark records the
# last edit position).
$p1 = "Divided up the main tasks we need to complete in order to start the implementation as soon as possible. We" + [char]0x2019 + "ve agreed to make a chat on facebook purely to upd"
$p2 = "ate each other on c"
$p3 = "hat documents are being worked on and when they are completed."

# Replace the whole paragraph's text, then re-insert it as three distinct
# runs (this mirrors Word splitting the run at the edited span).
$whole = $d.Range($rStart, $rEnd)
$whole.Text = ""

$seg1 = $d.Range($rStart, $rStart)
$seg1.InsertAfter($p1)

$seg2 = $d.Range($rStart + $p1.Length, $rStart + $p1.Length)
$seg2.InsertAfter($p2)

$seg3 = $d.Range($rStart + $p1.Length + $p2.Length, $rStart + $p1.Length + $p2.Length)
$seg3.InsertAfter($p3)

# Move the "_GoBack" bookmark (last-edit marker) to sit right after the
# corrected character, between the 2nd and 3rd runs.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmPos = $rStart + $p1.Length + $p2.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
